$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row at 10 (week of 2021-08-19)
$ws.Rows(10).Insert()
$ws.Cells.Item(10,1).Value = 2
$ws.Cells.Item(10,2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(10,3).Value = 'Coquimbo'
$ws.Cells.Item(10,4).Value = 44427
$ws.Cells.Item(10,5).Value = 4
$ws.Cells.Item(10,6).Value = 100112022
$ws.Cells.Item(10,7).Value = 'Arveja Verde'
$ws.Cells.Item(10,8).Value = 'Perfection'
$ws.Cells.Item(10,9).Value = 'Primera'
$ws.Cells.Item(10,10).Value = 300
$ws.Cells.Item(10,11).Value = 28000
$ws.Cells.Item(10,12).Value = 30000
$ws.Cells.Item(10,13).Value = 29000
$ws.Cells.Item(10,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(10,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(10,16).Value = 1160
$ws.Cells.Item(10,17).Value = 25
$ws.Cells.Item(10,18).Value = 'Hortaliza'

# Insert new row at 13 (week of 2021-08-18)
$ws.Rows(13).Insert()
$ws.Cells.Item(13,1).Value = 2
$ws.Cells.Item(13,2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(13,3).Value = 'Coquimbo'
$ws.Cells.Item(13,4).Value = 44426
$ws.Cells.Item(13,5).Value = 4
$ws.Cells.Item(13,6).Value = 100112022
$ws.Cells.Item(13,7).Value = 'Arveja Verde'
$ws.Cells.Item(13,8).Value = 'Perfection'
$ws.Cells.Item(13,9).Value = 'Primera'
$ws.Cells.Item(13,10).Value = 400
$ws.Cells.Item(13,11).Value = 28000
$ws.Cells.Item(13,12).Value = 30000
$ws.Cells.Item(13,13).Value = 29000
$ws.Cells.Item(13,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(13,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(13,16).Value = 1160
$ws.Cells.Item(13,17).Value = 25
$ws.Cells.Item(13,18).Value = 'Hortaliza'

# Insert 2 new rows at 21 (weeks of 2021-08-11 and 2021-08-12)
$ws.Rows(21).Insert()
$ws.Rows(21).Insert()
$ws.Cells.Item(21,1).Value = 2
$ws.Cells.Item(21,2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(21,3).Value = 'Coquimbo'
$ws.Cells.Item(21,4).Value = 44419
$ws.Cells.Item(21,5).Value = 4
$ws.Cells.Item(21,6).Value = 100112022
$ws.Cells.Item(21,7).Value = 'Arveja Verde'
$ws.Cells.Item(21,8).Value = 'Perfection'
$ws.Cells.Item(21,9).Value = 'Primera'
$ws.Cells.Item(21,10).Value = 600
$ws.Cells.Item(21,11).Value = 27000
$ws.Cells.Item(21,12).Value = 29000
$ws.Cells.Item(21,13).Value = 28000
$ws.Cells.Item(21,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(21,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(21,16).Value = 1120
$ws.Cells.Item(21,17).Value = 25
$ws.Cells.Item(21,18).Value = 'Hortaliza'
$ws.Cells.Item(22,1).Value = 2
$ws.Cells.Item(22,2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(22,3).Value = 'Coquimbo'
$ws.Cells.Item(22,4).Value = 44420
$ws.Cells.Item(22,5).Value = 4
$ws.Cells.Item(22,6).Value = 100112022
$ws.Cells.Item(22,7).Value = 'Arveja Verde'
$ws.Cells.Item(22,8).Value = 'Perfection'
$ws.Cells.Item(22,9).Value = 'Primera'
$ws.Cells.Item(22,10).Value = 700
$ws.Cells.Item(22,11).Value = 27000
$ws.Cells.Item(22,12).Value = 29000
$ws.Cells.Item(22,13).Value = 28000
$ws.Cells.Item(22,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(22,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(22,16).Value = 1120
$ws.Cells.Item(22,17).Value = 25
$ws.Cells.Item(22,18).Value = 'Hortaliza'
